$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '28.053.92'
$ws.Range('E2').Value = '  +3.29%  '
$ws.Range('D3').Value = '1.573.73'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -1.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.67'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.491'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  -1.15%  '
$ws.Range('E8').Value = '  +5.47%  '
$ws.Range('E9').Value = '  +0.55%  '
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('E11').Value = '  +1.69%  '
$ws.Range('D12').Value = '1.797.19'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = '1.574.16'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.76'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').Value = '28.028.47'
$ws.Range('E16').Value = '  +3.35%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.45'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '228.79'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +6.05%  '
$ws.Range('D19').Value = '0.0₃0705'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.44'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.32'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.93'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.31'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.21'
$ws.Range('D26').ClearFormats()
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.13'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.73%  '
$ws.Range('D34').Value = '1.416.51'
$ws.Range('E34').Value = '  -2.62%  '
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('E36').Value = '  -4.23%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.32'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.541'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.45'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.58%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.806'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.997'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.25%  '
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('E44').Value = '  -2.68%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.82'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.16%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '63.76'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('D47').Value = '1.710.34'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('E49').Value = '  +2.59%  '
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('E51').Value = '  -1.90%  '
